# "Loan RBI, Variable Instalments"
#
# The user switches to the "Repayment schedule" sheet and inserts a new
# (blank) column before column N, shifting the existing "Late",
# "heading"/"Outstanding(orig.)" and "Outstanding" columns one place to
# the right (N->O, O->P, P->Q). The new column is left blank/unused.
# Finally the selection on "Repayment schedule" lands on S11, which
# becomes the active sheet/tab (replacing "Transactions" as the
# previously-active tab).

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet (was "Transactions" before).
$wsRepayment.Activate() | Out-Null

# Insert a new blank column before column N (14th column); everything
# from N onward (N, O, P) shifts right by one (-> O, P, Q).
$wsRepayment.Columns("N:N").Insert() | Out-Null

# Excel copies the format of the column to the left (M) onto the freshly
# inserted column; reproduce that width here.
$wsRepayment.Columns(14).ColumnWidth = $wsRepayment.Columns(13).ColumnWidth

# Reflect the final on-sheet selection left behind after the edit.
$wsRepayment.Range("S11").Select() | Out-Null
